# Applies the "registros.xlsx" edit described in the commit:
#   "solucionando cantidad de caracteres que se pueden escribir"
# -> removes two stray empty placeholder cells in row 101 (Q101, Z101)
# -> appends 8 new data rows (102-109) to the "Obras en general" sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Obras en general")

function Set-TextCell {
    param($sheet, [int]$row, [int]$col, [string]$text)

    $cell = $sheet.Cells.Item($row, $col)
    # Force text storage so dates / numeric-looking strings ("12/06/2025",
    # "18794", "1911", ...) are kept as literal text instead of being
    # auto-converted to a date serial or a number. Setting NumberFormat
    # stamps a style index on the cell (and some columns also carry a
    # column-level style, e.g. F/G/L/M/N/O/P/Q); the source file never
    # stamps an explicit style on these data cells, so reset back to the
    # plain "Normal" style once the text value has been written.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Row 101 clean-up: two empty placeholder cells (Q101, Z101) are removed
# ---------------------------------------------------------------------
$ws.Cells.Item(101, 17).ClearContents()   # Q101
$ws.Cells.Item(101, 26).ClearContents()   # Z101

# ---------------------------------------------------------------------
# New rows 102-109
# ---------------------------------------------------------------------
$newRows = @(
    @{
        row = 102
        A = "12/06/2025"; B = "MMO"; C = "Digital"
        E = "Obra nueva"; F = "JUANJO ROMERO"
        G = "ORIANA MENDOZA / SANTIAGOA MENDOZA / LAUTARO MENDOZA / WALTER MENDOZA / MARIELA DIAZ MENDOZA"
        H = "VILLA SARITA"; I = "12344/R/25"; J = "1911"; K = "18794"
        R = "No pagado"; S = "No pagado"
        X = "c:\Users\Admin\Desktop\Gestion-expedientes-cpim\trabajos\OBRA NUEVA\JUANJO ROMERO\ORIANA MENDOZA _ SANTIAGOA MENDOZA _ LAUTARO MENDOZA _ WALTER MENDOZA _ MARIELA DIAZ MENDOZA"
        Y = "3764251817"
    },
    @{
        row = 103
        A = "13/06/2025"; B = "Ingeniero"; C = "Físico"; D = "2"
        E = "Registración"; F = "DE JESUS SANTIAGO"
        G = "MENDOZA ORIANA / MENDOZA SANTIAGO / MENDOZA MARIELA / MENDOZA WALTER / MENDOZAAA"
        H = "MDISDIA"; I = "1818/F/181"; K = "118"
        R = "No pagado"; S = "No pagado"
    },
    @{
        row = 104
        A = "13/06/2025"; B = "Ingeniero"; C = "Físico"; D = "2"
        E = "Obra nueva"; F = "SANTIAGOOOOO"
        G = "ORIANA MENDOZA / MENDOZA NOSE / MENDOZA LAUTADROOO"
        H = "GVGERGE"; I = "18/R/G21"; K = "651"
        R = "No pagado"; S = "No pagado"
    },
    @{
        row = 105
        A = "13/06/2025"; B = "Licenciado"; C = "Físico"; D = "2"
        E = "Obra nueva"; F = "DIFI"
        G = "FFRFERIFMERIFFERFERFGREG"
        H = "GERGERGERG"; I = "188/GR/185"; K = "1789"
        R = "No pagado"; S = "No pagado"
    },
    @{
        row = 106
        A = "13/06/2025"; B = "Licenciado"; C = "Físico"; D = "1"
        E = "Obra nueva"; F = "LALALALALA"
        G = "LALALA / LALALA"
        H = "RGERGER"; I = "GERG/8GER"; K = "5145"
        R = "No pagado"; S = "No pagado"
    },
    @{
        row = 107
        A = "13/06/2025"; B = "MMO"; C = "Físico"; D = "1"
        E = "Obra nueva"; F = "LALALAL"
        G = "LALALA - LALALA - LALAL"
        H = "FEWFGWEG"; I = "GERGV/185"; K = "5185"
        R = "No pagado"; S = "No pagado"
    },
    @{
        row = 108
        A = "13/06/2025"; B = "Licenciado"; C = "Físico"; D = "1"
        E = "Obra nueva"; F = "LALALALALITA"
        G = "LALALAL - LALALALA - LALALALA - LALALALAL - LALALA"
        H = "REGERG"; I = "GERGBVER/8GERERG"; K = "8418"
        R = "No pagado"; S = "No pagado"
    },
    @{
        row = 109
        A = "13/06/2025"; B = "Licenciado"; C = "Físico"; D = "178"
        E = "Obra nueva"; F = "benitez lucia ines"
        G = "CARNICERIA - LA GRANKA - SIEMPRE PRECIOS BBAJOS"
        H = "FGERWGRE"; I = "EWFWE/8FEW/"; K = "FWEFWE"
        R = "No pagado"; S = "No pagado"
    }
)

# Column letter -> index map for the fields used above
$colIndex = @{
    A=1; B=2; C=3; D=4; E=5; F=6; G=7; H=8; I=9; J=10; K=11; L=12; M=13
    N=14; O=15; P=16; Q=17; R=18; S=19; T=20; U=21; V=22; W=23; X=24; Y=25
    Z=26; AA=27
}

foreach ($rdata in $newRows) {
    $rowNum = $rdata.row
    foreach ($key in $colIndex.Keys) {
        if ($rdata.ContainsKey($key)) {
            Set-TextCell $ws $rowNum $colIndex[$key] $rdata[$key]
        }
    }
}
